$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.387.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.50%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.957.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +3.36%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.53%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'328.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.85%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +0.40%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4642"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.20%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3935"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.74%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -0.52%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.07900"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.59%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'1.002"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.22%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'22.37"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +2.13%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.976.53"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +3.18%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'7.162"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.46%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'5.851"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +2.78%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.07143"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +3.00%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'88.55"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.57%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'1.006"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.53%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.000009946"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.40%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'17.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.66%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'1.004"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.40%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'29.439.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.64%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'5.510"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +3.96%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'11.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +2.51%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.154.83"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.15%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.125"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +3.14%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'158.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.50%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'19.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.48%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D30").Value = "'119.72"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +1.79%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'1.892"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -2.01%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.09378"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.18%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.8968"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.20%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'5.248"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.90%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.337"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.31%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'3.190"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -2.18%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.000003781"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +118.06%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.05825"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.94%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'1.175"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.37%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.02120"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +2.08%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +0.53%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'7.786"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.46%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.5756"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.24%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.1824"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +2.87%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'9.797"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.43%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'12.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +2.40%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.5377"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.37%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'2.212"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -4.09%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.890"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +2.34%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'2.614"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +3.45%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.06957"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.26%  "
$ws.Range("E51").Style = "Normal"
